$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.868.57"
$ws.Range("E2").Value = "  -1.09%  "
Set-TextValue $ws.Range("D3") "3.054.33"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws.Range("D5") "559.64"
$ws.Range("E5").Value = "  +0.21%  "
Set-TextValue $ws.Range("D6") "142.74"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +0.09%  "
Set-TextValue $ws.Range("D8") "3.053.63"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("E10").Value = "  +0.88%  "
Set-TextValue $ws.Range("D11") "6.27"
$ws.Range("E11").Value = "  -11.07%  "
Set-TextValue $ws.Range("D12") "0.491"
$ws.Range("E12").Value = "  +6.59%  "
$ws.Range("E13").Value = "  +1.66%  "
Set-TextValue $ws.Range("D14") "35.76"
$ws.Range("E14").Value = "  +1.88%  "
Set-TextValue $ws.Range("D15") "3.552.85"
$ws.Range("E15").Value = "  -1.25%  "
Set-TextValue $ws.Range("D16") "63.925.77"
$ws.Range("E16").Value = "  -0.86%  "
Set-TextValue $ws.Range("D17") "3.054.31"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("E19").Value = "  +1.34%  "
Set-TextValue $ws.Range("D20") "476.72"
$ws.Range("E20").Value = "  -1.49%  "
Set-TextValue $ws.Range("D21") "14.08"
$ws.Range("E21").Value = "  +2.54%  "
Set-TextValue $ws.Range("D22") "14.70"
$ws.Range("E22").Value = "  +11.20%  "
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  +0.01%  "
Set-TextValue $ws.Range("D25") "82.78"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  -0.06%  "
Set-TextValue $ws.Range("D27") "2.80"
$ws.Range("E27").Value = "  -0.49%  "
Set-TextValue $ws.Range("D28") "8.18"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("E29").Value = "  -1.09%  "
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.29%  "
Set-TextValue $ws.Range("D31") "26.31"
$ws.Range("E31").Value = "  +0.78%  "
Set-TextValue $ws.Range("D32") "1.15"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +1.85%  "
Set-TextValue $ws.Range("D36") "54.59"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  +1.21%  "
Set-TextValue $ws.Range("D38") "447.59"
$ws.Range("E38").Value = "  -3.36%  "
Set-TextValue $ws.Range("D39") "0.0816"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  +4.06%  "
Set-TextValue $ws.Range("D41") "3.018.88"
$ws.Range("E41").Value = "  -0.19%  "
Set-TextValue $ws.Range("D42") "0.118"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +4.36%  "
Set-TextValue $ws.Range("D45") "28.22"
$ws.Range("E45").Value = "  +0.23%  "
Set-TextValue $ws.Range("D46") "2.28"
$ws.Range("E46").Value = "  +9.64%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +1.10%  "
Set-TextValue $ws.Range("D49") "117.82"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("E51").Value = "  +1.51%  "
